# Lab Exam 03 grading workbook - "Changes to driver class 65-80"
#
# The grader went back through the "CustomerMappingDriver Class" and
# "Compilation errors" rows and rewrote a couple of grading comments,
# then docked additional points for a NullPointerException found while
# re-checking compilation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CustomerMappingDriver Class section (rows 29-30): replace the generic
# scanner/output comments with more specific ones about the driver.
$ws.Range("F29").Value = "(-16) for not attempting to write anything in driver"
$ws.Range("F30").Value = "(-4) for no output displayed due to NullPointerException"

# Compilation errors row (37): extra deduction + comment for the
# NullPointerException uncovered while compiling.
$ws.Range("E37").Value = -2.5
$ws.Range("F37").Value = "(-2.5) for getting NullPointerException"

# Leave the view scrolled/selected where the grader left off.
$ws.Range("F37").Select() | Out-Null
